# Scheduled market-data refresh: rewrite the price/profit columns (H:N)
# for each updated leve row across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 674.1429000000001
$ws.Range("I8").Value = 674.1429000000001
$ws.Range("K8").Value = 2022.4287
$ws.Range("M8").Value = -1883.4287
$ws.Range("H12").Value = 474.27274
$ws.Range("I12").Value = 91.833336
$ws.Range("J12").Value = 933.2
$ws.Range("K12").Value = 91.833336
$ws.Range("L12").Value = 933.2
$ws.Range("M12").Value = 78.166664
$ws.Range("N12").Value = -1273.2
$ws.Range("I43").Value = 1299
$ws.Range("J43").Value = 2311
$ws.Range("K43").Value = 1299
$ws.Range("L43").Value = 2311
$ws.Range("M43").Value = -1230
$ws.Range("N43").Value = -2449
$ws.Range("H74").Value = 6260.8
$ws.Range("I74").Value = 5479.9
$ws.Range("J74").Value = 7041.7
$ws.Range("K74").Value = 5479.9
$ws.Range("L74").Value = 7041.7
$ws.Range("M74").Value = -4543.9
$ws.Range("N74").Value = -8913.700000000001
$ws.Range("H77").Value = 6260.8
$ws.Range("I77").Value = 5479.9
$ws.Range("J77").Value = 7041.7
$ws.Range("K77").Value = 27399.5
$ws.Range("L77").Value = 35208.5
$ws.Range("M77").Value = -22719.5
$ws.Range("N77").Value = -44568.5
$ws.Range("H96").Value = 730.8889
$ws.Range("I96").Value = 446.66666
$ws.Range("K96").Value = 1339.99998
$ws.Range("M96").Value = 33.00001999999995
$ws.Range("H116").Value = 11480
$ws.Range("I116").Value = 18779
$ws.Range("J116").Value = 7100.6
$ws.Range("K116").Value = 18779
$ws.Range("L116").Value = 7100.6
$ws.Range("M116").Value = -15337
$ws.Range("N116").Value = -13984.6
$ws.Range("H136").Value = 54285.57
$ws.Range("J136").Value = 54285.57
$ws.Range("L136").Value = 54285.57
$ws.Range("N136").Value = -64485.57
$ws.Range("H138").Value = 4027.543
$ws.Range("I138").Value = 3265.44
$ws.Range("J138").Value = 5932.8
$ws.Range("K138").Value = 9796.32
$ws.Range("L138").Value = 17798.4
$ws.Range("M138").Value = -4656.32
$ws.Range("N138").Value = -28078.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1183.2307
$ws.Range("I2").Value = 670.6667
$ws.Range("J2").Value = 1622.5714
$ws.Range("K2").Value = 670.6667
$ws.Range("L2").Value = 1622.5714
$ws.Range("M2").Value = -557.6667
$ws.Range("N2").Value = -1848.5714
$ws.Range("H32").Value = 2073.0923
$ws.Range("I32").Value = 1928.9678
$ws.Range("K32").Value = 1928.9678
$ws.Range("M32").Value = -1641.9678
$ws.Range("H116").Value = 1183.2307
$ws.Range("I116").Value = 670.6667
$ws.Range("J116").Value = 1622.5714
$ws.Range("K116").Value = 670.6667
$ws.Range("L116").Value = 1622.5714
$ws.Range("M116").Value = 1623.3333
$ws.Range("N116").Value = -6210.5714
$ws.Range("H132").Value = 5321.7896
$ws.Range("J132").Value = 9337.666999999999
$ws.Range("L132").Value = 28013.001
$ws.Range("N132").Value = -33073.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1183.2307
$ws.Range("I3").Value = 670.6667
$ws.Range("J3").Value = 1622.5714
$ws.Range("K3").Value = 670.6667
$ws.Range("L3").Value = 1622.5714
$ws.Range("M3").Value = -556.6667
$ws.Range("N3").Value = -1850.5714
$ws.Range("H132").Value = 90999
$ws.Range("J132").Value = 90999
$ws.Range("L132").Value = 90999
$ws.Range("N132").Value = -101119
$ws.Range("H133").Value = 53994.75
$ws.Range("I133").Value = 35000
$ws.Range("J133").Value = 60326.332
$ws.Range("K133").Value = 35000
$ws.Range("L133").Value = 60326.332
$ws.Range("M133").Value = -29940
$ws.Range("N133").Value = -70446.33199999999
$ws.Range("H139").Value = 161784.33
$ws.Range("J139").Value = 183999.4
$ws.Range("L139").Value = 183999.4
$ws.Range("N139").Value = -194279.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 17198.4
$ws.Range("J17").Value = 17198.4
$ws.Range("L17").Value = 17198.4
$ws.Range("N17").Value = -17546.4
$ws.Range("H58").Value = 1274.5
$ws.Range("I58").Value = 1274.5
$ws.Range("K58").Value = 1274.5
$ws.Range("M58").Value = -1071.5
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
$ws.Range("H132").Value = 4340
$ws.Range("I132").Value = 4340
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13020
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10490
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1274.5
$ws.Range("I136").Value = 1274.5
$ws.Range("K136").Value = 3823.5
$ws.Range("M136").Value = -1273.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 522
$ws.Range("J23").Value = 589
$ws.Range("L23").Value = 1767
$ws.Range("N23").Value = -2237
$ws.Range("H50").Value = 1965.625
$ws.Range("J50").Value = 900
$ws.Range("L50").Value = 2700
$ws.Range("N50").Value = -3662
$ws.Range("H53").Value = 1965.625
$ws.Range("J53").Value = 900
$ws.Range("L53").Value = 2700
$ws.Range("N53").Value = -3662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 14000
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 15000
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = -9748
$ws.Range("N33").Value = -15504
$ws.Range("H40").Value = 12666.333
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10302
$ws.Range("H80").Value = 2557
$ws.Range("J80").Value = 2499.75
$ws.Range("L80").Value = 2499.75
$ws.Range("N80").Value = -4495.75
$ws.Range("H83").Value = 2557
$ws.Range("J83").Value = 2499.75
$ws.Range("L83").Value = 12498.75
$ws.Range("N83").Value = -22482.75
$ws.Range("H107").Value = 442.75
$ws.Range("I107").Value = 242.35715
$ws.Range("J107").Value = 723.3
$ws.Range("K107").Value = 242.35715
$ws.Range("L107").Value = 723.3
$ws.Range("M107").Value = 1677.64285
$ws.Range("N107").Value = -4563.3
$ws.Range("H122").Value = 2338.842
$ws.Range("I122").Value = 2303.1538
$ws.Range("J122").Value = 2416.1667
$ws.Range("K122").Value = 6909.4614
$ws.Range("L122").Value = 7248.500100000001
$ws.Range("M122").Value = -4459.4614
$ws.Range("N122").Value = -12148.5001
$ws.Range("H126").Value = 6644.1333
$ws.Range("I126").Value = 5158.1113
$ws.Range("K126").Value = 15474.3339
$ws.Range("M126").Value = -13004.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3184.7827
$ws.Range("I132").Value = 3992.6
$ws.Range("J132").Value = 2960.389
$ws.Range("K132").Value = 11977.8
$ws.Range("L132").Value = 8881.167000000001
$ws.Range("M132").Value = -9447.799999999999
$ws.Range("N132").Value = -13941.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 59999.5
$ws.Range("J70").Value = 59999.5
$ws.Range("L70").Value = 59999.5
$ws.Range("N70").Value = -60629.5
$ws.Range("H73").Value = 59999.5
$ws.Range("J73").Value = 59999.5
$ws.Range("L73").Value = 59999.5
$ws.Range("N73").Value = -62183.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
